$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ზუგდიდი")

$ws.Range("E4").Value = 5189
$ws.Range("F4").Value = 5283
$ws.Range("G4").Value = 5280
$ws.Range("H4").Value = 5231
